$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the daily report by appending 6 more days (rows 27-32), following
# the same layout/format as the existing rows. First clone the formatting of
# the last existing row down into the new rows, then fill in the data.
$ws.Range("A26:Q26").Copy()
$ws.Range("A27:A32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=27; A=26; B=45744; Day="Fri"; D=9353; E=9353; L=0; O=87249.1; P=84105.3 },
    @{ Row=28; A=27; B=45745; Day="Sat"; D=9353; E=9400; L=0; O=84105.3; P=82428.8 },
    @{ Row=29; A=28; B=45746; Day="Sun"; D=9400; E=9450; L=0; O=82428.8; P=82745.1 },
    @{ Row=30; A=29; B=45747; Day="Mon"; D=9450; E=9452; L=0; O=82745.1; P=83448.1 },
    @{ Row=31; A=30; B=45748; Day="Tue"; D=9452; E=9454; L=0; O=83448.1; P=84138.9 },
    @{ Row=32; A=31; B=45749; Day="Wed"; D=9454; E=9456; L=0; O=84138.9; P=85827.9 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.Day
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    $ws.Cells.Item($row, 6).Formula = "=E$row-D$row"
    $ws.Cells.Item($row, 7).Formula = "=(E$row-`$D`$2)/A$row"
    $ws.Cells.Item($row, 8).Formula = "=(E$row/D$row-1)*100"
    $ws.Cells.Item($row, 9).Formula = "=(POWER((E$row/`$D`$3),1/A$row)-1)*100"
    $ws.Cells.Item($row, 10).Formula = "=J$prevRow*1.013"
    $ws.Cells.Item($row, 11).Formula = "=E$row-J$row"

    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Formula = "=L$row+E$row"
    $ws.Cells.Item($row, 14).Formula = "=E$row/`$D`$2*100"

    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Formula = "=P$row/`$O`$2*100"
}

$ws.Range("M23").Select()
